$wb = $excel.ActiveWorkbook

# --- Sheet 1: Book Inventory ---
$ws1 = $wb.Worksheets.Item(1)

# Save old row2 ISBN (C2) into a scratch cell so we can move it to C3 later.
$ws1.Range("C2").Copy()
$ws1.Range("G1").PasteSpecial(-4163)

# Move old row3 (A3:C3) "The Hunger Games" data up into row2.
$ws1.Range("A3:C3").Copy()
$ws1.Range("A2").PasteSpecial(-4163)

# New quantities for "The Hunger Games" row (now row 2).
$ws1.Range("D2").Value = 2
$ws1.Range("E2").Value = 2

# Row 3 becomes the "Wild Things Storytelling Kit" companion item.
$ws1.Range("A3").Value = "Wild Things Storytelling Kit"
$ws1.Range("B3").Value = "story and pictures by Maurice Sendak"
$ws1.Range("G1").Copy()
$ws1.Range("C3").PasteSpecial(-4163)

# Clear scratch cell.
$ws1.Range("G1").Clear()

# --- Sheet 2: Check Out-In (check the book back in) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(2).Delete()
